# Child Obligation Suite.xlsx - 24-Feb-2017 edit
# - Update a few Runmode/Results cells on the three sheets
# - Move the active sheet / selection from "Test Cases" to "ChildObligationWFNonWF"

$wb = $excel.ActiveWorkbook

$wsTestCases   = $wb.Worksheets.Item(1)   # "Test Cases"
$wsChildDNO    = $wb.Worksheets.Item(2)   # "ChildDNOUpdate"
$wsChildWF     = $wb.Worksheets.Item(3)   # "ChildObligationWFNonWF"

# --- Test Cases sheet: ChildObligationWFNonWF run now passes, ChildDNOUpdate run now fails ---
$wsTestCases.Range("C2").Value = "Y"
$wsTestCases.Range("D2").Value = "FAIL"
$wsTestCases.Range("D3").Value = "FAIL"

# --- ChildDNOUpdate sheet: result flips to FAIL ---
$wsChildDNO.Range("G2").Value = "FAIL"

# --- ChildObligationWFNonWF sheet: OnHold no longer run, overall result is now FAIL ---
$wsChildWF.Range("G2").Value = ""
$wsChildWF.Range("K2").Value = "FAIL"

# --- Update selections on "Test Cases" (no longer the active tab) and on the new active tab ---
$wsTestCases.Activate()
$wsTestCases.Range("C2").Select()

$wsChildWF.Activate()
$wsChildWF.Range("C2").Select()
